$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set the Performer for the "Oriented Graph class realization" task (row 3)
$ws.Range("C3").Value = "Kochetov Pavel"

# Update the active selection on the sheet to C4
$ws.Range("C4").Select()
